$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected (legacy password hash "D382"); unprotect to edit,
# then restore protection afterwards.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclaimer text (A10):
# 2021-05-18 -> 2021-05-19
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-19 for illustrative purposes only and are subject to change."

# Re-fit row 10's height so the multi-line text doesn't leave a stray
# custom row-height behind (matches original unset/auto height).
$ws.Rows.Item(10).AutoFit()

# Update the Weight (D) / Percent Change (E) figures for rows 2-7.
$ws.Range("D2").Value = 0.2429035270437203
$ws.Range("E2").Value = 0.000881171616710974

$ws.Range("D3").Value = 0.502740691093806
$ws.Range("E3").Value = -0.005540897097625286

$ws.Range("D4").Value = 0.09486925400713392
$ws.Range("E4").Value = 0.0004086587092859695

$ws.Range("D5").Value = 0.1028473184103282
$ws.Range("E5").Value = -0.009513454743136696

$ws.Range("D6").Value = 0.05663920944501157
$ws.Range("E6").Value = -0.007537688442211032

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = -0.003938187619175748

# Restore sheet protection with the original password.
$ws.Protect("D382")
